$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 308.91306
$ws.Range("I55").Value = 202.75
$ws.Range("J55").Value = 365.53333
$ws.Range("K55").Value = 202.75
$ws.Range("L55").Value = 365.53333
$ws.Range("M55").Value = 11.25
$ws.Range("N55").Value = -793.53333
$ws.Range("H70").Value = 1966.6666
$ws.Range("I70").Value = 1950
$ws.Range("J70").Value = 2000
$ws.Range("K70").Value = 5850
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -5580
$ws.Range("N70").Value = -6540
$ws.Range("H73").Value = 1966.6666
$ws.Range("I73").Value = 1950
$ws.Range("J73").Value = 2000
$ws.Range("K73").Value = 5850
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -4914
$ws.Range("N73").Value = -7872
$ws.Range("H86").Value = 1941
$ws.Range("I86").Value = 1926.25
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1926.25
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -803.25
$ws.Range("N86").Value = -4246
$ws.Range("H89").Value = 1941
$ws.Range("I89").Value = 1926.25
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 9631.25
$ws.Range("L89").Value = 2000
$ws.Range("M89").Value = -4015.25
$ws.Range("N89").Value = -21232
$ws.Range("H112").Value = 1261.5363
$ws.Range("J112").Value = 1261.5363
$ws.Range("L112").Value = 3784.6089
$ws.Range("N112").Value = -6000.6089
$ws.Range("H127").Value = 1164.1111
$ws.Range("I127").Value = 955
$ws.Range("J127").Value = 1425.5
$ws.Range("K127").Value = 2865
$ws.Range("L127").Value = 4276.5
$ws.Range("M127").Value = 2095
$ws.Range("N127").Value = -14196.5
$ws.Range("H137").Value = 838459.4
$ws.Range("I137").Value = 1538967.6
$ws.Range("K137").Value = 4616902.800000001
$ws.Range("M137").Value = -4614352.800000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4714.5938
$ws.Range("I32").Value = 5083.196
$ws.Range("K32").Value = 5083.196
$ws.Range("M32").Value = -4796.196
$ws.Range("H74").Value = 411636.44
$ws.Range("I74").Value = 592942.4
$ws.Range("J74").Value = 3698
$ws.Range("K74").Value = 592942.4
$ws.Range("L74").Value = 3698
$ws.Range("M74").Value = -592068.4
$ws.Range("N74").Value = -5446
$ws.Range("H77").Value = 411636.44
$ws.Range("I77").Value = 592942.4
$ws.Range("J77").Value = 3698
$ws.Range("K77").Value = 2964712
$ws.Range("L77").Value = 18490
$ws.Range("M77").Value = -2960344
$ws.Range("N77").Value = -27226
$ws.Range("H102").Value = 2300.7273
$ws.Range("I102").Value = 2300.7273
$ws.Range("K102").Value = 2300.7273
$ws.Range("M102").Value = -678.7273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 4777.5
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 4777.5
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 4777.5
$ws.Range("M5").ClearContents()
$ws.Range("N5").Value = -5003.5
$ws.Range("H7").Value = 7875714.5
$ws.Range("I7").Value = 15000000
$ws.Range("J7").Value = 6688334
$ws.Range("K7").Value = 15000000
$ws.Range("L7").Value = 6688334
$ws.Range("M7").Value = -14999887
$ws.Range("N7").Value = -6688560

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 248007.89
$ws.Range("I31").Value = 564022.94
$ws.Range("J31").Value = 3351.0645
$ws.Range("K31").Value = 564022.94
$ws.Range("L31").Value = 3351.0645
$ws.Range("M31").Value = -563727.94
$ws.Range("N31").Value = -3941.0645
$ws.Range("H34").Value = 248007.89
$ws.Range("I34").Value = 564022.94
$ws.Range("J34").Value = 3351.0645
$ws.Range("K34").Value = 564022.94
$ws.Range("L34").Value = 3351.0645
$ws.Range("M34").Value = -563820.94
$ws.Range("N34").Value = -3755.0645
$ws.Range("H38").Value = 118099.91
$ws.Range("I38").Value = 1000000
$ws.Range("J38").Value = 29909.9
$ws.Range("K38").Value = 1000000
$ws.Range("L38").Value = 29909.9
$ws.Range("M38").Value = -999623
$ws.Range("N38").Value = -30663.9
$ws.Range("H46").Value = 118099.91
$ws.Range("I46").Value = 1000000
$ws.Range("J46").Value = 29909.9
$ws.Range("K46").Value = 1000000
$ws.Range("L46").Value = 29909.9
$ws.Range("M46").Value = -999789
$ws.Range("N46").Value = -30331.9
$ws.Range("H52").Value = 57475
$ws.Range("J52").Value = 57475
$ws.Range("L52").Value = 57475
$ws.Range("N52").Value = -58063
$ws.Range("H60").Value = 24588
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 24588
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 24588
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -25610

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2330.6667
$ws.Range("I5").Value = 592.6667
$ws.Range("J5").Value = 2910
$ws.Range("K5").Value = 1778.0001
$ws.Range("L5").Value = 8730
$ws.Range("M5").Value = -1666.0001
$ws.Range("N5").Value = -8954
$ws.Range("H22").Value = 294.83334
$ws.Range("J22").Value = 290
$ws.Range("L22").Value = 870
$ws.Range("N22").Value = -1208
$ws.Range("H27").Value = 294.83334
$ws.Range("J27").Value = 290
$ws.Range("L27").Value = 870
$ws.Range("N27").Value = -1074
$ws.Range("H135").Value = 2330.6667
$ws.Range("I135").Value = 592.6667
$ws.Range("J135").Value = 2910
$ws.Range("K135").Value = 5334.0003
$ws.Range("L135").Value = 26190
$ws.Range("M135").Value = -2799.0003
$ws.Range("N135").Value = -31260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 3266.6667
$ws.Range("J19").Value = 12950
$ws.Range("L19").Value = 12950
$ws.Range("N19").Value = -13526
$ws.Range("H107").Value = 6945235
$ws.Range("I107").Value = 565.5714
$ws.Range("J107").Value = 12346644
$ws.Range("K107").Value = 565.5714
$ws.Range("L107").Value = 12346644
$ws.Range("M107").Value = 1354.4286
$ws.Range("N107").Value = -12350484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5359.1904
$ws.Range("I40").Value = 4742.533
$ws.Range("J40").Value = 6900.8335
$ws.Range("K40").Value = 4742.533
$ws.Range("L40").Value = 6900.8335
$ws.Range("M40").Value = -4606.533
$ws.Range("N40").Value = -7172.8335
$ws.Range("H100").Value = 4970.4
$ws.Range("I100").Value = 2140
$ws.Range("J100").Value = 7800.8
$ws.Range("K100").Value = 2140
$ws.Range("L100").Value = 7800.8
$ws.Range("M100").Value = -1599
$ws.Range("N100").Value = -8882.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5459.189
$ws.Range("I136").Value = 5231.409
$ws.Range("J136").Value = 5793.2666
$ws.Range("K136").Value = 15694.227
$ws.Range("L136").Value = 17379.7998
$ws.Range("M136").Value = -13144.227
$ws.Range("N136").Value = -22479.7998
